# Rubric-Sprint 3v1.xlsx: "Batch process to send emails to rentals overdue"
# (row 11) previously had no "Possible points" score in column C; fill it
# in with 8 (matching the "Points" value already in column D), which in
# turn changes the dependent SUM()/percentage formulas further down the
# sheet (E7, E25, E27, F28) -- those recalculate automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = 8

# Make sure every formula that depends on C11 (directly or transitively)
# is recalculated before the workbook is saved.
$excel.Calculate()

# Reflect where the editor left the cursor/scroll position.
$ws.Range("B12").Select()
try {
    $excel.ActiveWindow.ScrollRow = 2
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # Older/limited hosts may not expose ActiveWindow scroll properties;
    # the selection above is the important, persisted part.
}
